$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 101: Petty King / Klein-Koning / 2024-05-28
$ws.Range("A101").Value = "Petty King"
$ws.Range("B101").Value = "Klein-Koning"
$ws.Range("D101").Value = "2024-05-28"
$ws.Range("D101").NumberFormat = "@"
$ws.Range("D101").HorizontalAlignment = -4108

# Row 71: add note in column E about "Heldenmoed"
$ws.Range("E71").Value = "Edit: Heb ik vertaald naar Heldenmoed"

# New row 102: Court Physician / Hofarts / 2024-05-28
$ws.Range("A102").Value = "Court Physician"
$ws.Range("B102").Value = "Hofarts"
$ws.Range("D102").Value = "2024-05-28"
$ws.Range("D102").NumberFormat = "@"
$ws.Range("D102").HorizontalAlignment = -4108

# New row 103: Piety / Vroomheid / 2024-05-28
$ws.Range("A103").Value = "Piety"
$ws.Range("B103").Value = "Vroomheid"
$ws.Range("D103").Value = "2024-05-28"
$ws.Range("D103").NumberFormat = "@"
$ws.Range("D103").HorizontalAlignment = -4108

# New row 104: Level of Faith / Niveau van Toewijding / 2024-05-28
$ws.Range("A104").Value = "Level of Faith"
$ws.Range("B104").Value = "Niveau van Toewijding"
$ws.Range("D104").Value = "2024-05-28"
$ws.Range("D104").NumberFormat = "@"
$ws.Range("D104").HorizontalAlignment = -4108

# Update view state
$ws.Range("C103").Select()
$excel.ActiveWindow.ScrollRow = 76
